$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '237.49'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '21.64'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.465'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05647'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.490'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.353'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.071'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.7939'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1397'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07351'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03201'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09240'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001672'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.263'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04781'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005745'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006207'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.005107'
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001052'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001501'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.878'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.198'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004013'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04117'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006939'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1042'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003012'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009909'
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005443'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000751'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6758'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.03725'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002102'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01011'
